# Weekly Fruta/Hortaliza update: insert a new weekly record as row 7
# (Macroferia Regional de Talca, Espárragos) and push the existing
# historical rows (old rows 7-14) down to rows 8-15.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 7; this shifts the old
# rows 7-14 down to 8-15, matching the rest of the diff exactly.
$ws.Rows.Item(7).Insert()

# Populate the newly inserted row 7 with the new weekly record.
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "Macroferia Regional de Talca"
$ws.Range("C7").Value = "Maule"
$ws.Range("D7").Value = 44460
$ws.Range("E7").Value = 7
$ws.Range("F7").Value = 300000000
$ws.Range("G7").Value = "Espárragos"
$ws.Range("H7").Value = "Verde"
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 2000
$ws.Range("K7").Value = 2000
$ws.Range("L7").Value = 2000
$ws.Range("M7").Value = 2000
$ws.Range("N7").Value = "$/kilo"
$ws.Range("O7").Value = "Provincia de Linares"
$ws.Range("P7").Value = 2000
$ws.Range("Q7").Value = 1
$ws.Range("R7").Value = "Hortaliza"
